$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 200.5
$ws.Range("I5").Value = 100.666664
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 100.666664
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 14.333336
$ws.Range("N5").Value = -730
$ws.Range("H17").Value = 742838.0600000001
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 742838.0600000001
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2228514.18
$ws.Range("N17").Value = -2228850.18
$ws.Range("H33").Value = 4545866.5
$ws.Range("I33").Value = 4762310.5
$ws.Range("J33").Value = 540
$ws.Range("K33").Value = 4762310.5
$ws.Range("L33").Value = 540
$ws.Range("M33").Value = -4762081.5
$ws.Range("H37").Value = 94.25
$ws.Range("I37").Value = 94.25
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 282.75
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -156.75
$ws.Range("H86").Value = 5698.4
$ws.Range("I86").Value = 4750
$ws.Range("J86").Value = 6330.6665
$ws.Range("K86").Value = 4750
$ws.Range("L86").Value = 6330.6665
$ws.Range("M86").Value = -3627
$ws.Range("N86").Value = -8576.666499999999
$ws.Range("H89").Value = 5698.4
$ws.Range("I89").Value = 4750
$ws.Range("J89").Value = 6330.6665
$ws.Range("K89").Value = 23750
$ws.Range("L89").Value = 31653.3325
$ws.Range("M89").Value = -18134
$ws.Range("N89").Value = -42885.3325
$ws.Range("H92").Value = 542.2941
$ws.Range("I92").Value = 407.6154
$ws.Range("J92").Value = 980
$ws.Range("K92").Value = 407.6154
$ws.Range("L92").Value = 980
$ws.Range("M92").Value = 840.3846
$ws.Range("H106").Value = 3931.8
$ws.Range("I106").Value = 3931.8
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3931.8
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -3300.8
$ws.Range("H112").Value = 57447.332
$ws.Range("I112").Value = 1664.5
$ws.Range("J112").Value = 64420.188
$ws.Range("K112").Value = 4993.5
$ws.Range("L112").Value = 193260.564
$ws.Range("M112").Value = -3885.5
$ws.Range("N112").Value = -195476.564
$ws.Range("H137").Value = 1784.6666
$ws.Range("I137").Value = 1555
$ws.Range("J137").Value = 1899.5
$ws.Range("K137").Value = 4665
$ws.Range("L137").Value = 5698.5
$ws.Range("M137").Value = -2115
$ws.Range("N137").Value = -10798.5
$ws.Range("H138").Value = 3639.1912
$ws.Range("I138").Value = 2095.1304
$ws.Range("J138").Value = 4428.378
$ws.Range("K138").Value = 6285.3912
$ws.Range("L138").Value = 13285.134
$ws.Range("M138").Value = -1145.3912
$ws.Range("N138").Value = -23565.134
$ws.Range("H140").Value = 109998.2
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 109998.2
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 109998.2
$ws.Range("N140").Value = -120358.2
$ws.Range("H141").Value = 1325.2727
$ws.Range("I141").Value = 1325.2727
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3975.8181
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 1204.1819

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4435.205
$ws.Range("I32").Value = 4291.3945
$ws.Range("J32").Value = 9900
$ws.Range("K32").Value = 4291.3945
$ws.Range("L32").Value = 9900
$ws.Range("M32").Value = -4004.3945
$ws.Range("N32").Value = -10474
$ws.Range("H122").Value = 4135.6665
$ws.Range("I122").Value = 4103.4287
$ws.Range("J122").Value = 4248.5
$ws.Range("K122").Value = 12310.2861
$ws.Range("L122").Value = 12745.5
$ws.Range("M122").Value = -9860.286100000001
$ws.Range("N122").Value = -17645.5
$ws.Range("H132").Value = 3671.8215
$ws.Range("I132").Value = 3309
$ws.Range("J132").Value = 5002.1665
$ws.Range("K132").Value = 9927
$ws.Range("L132").Value = 15006.4995
$ws.Range("M132").Value = -7397

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2070.8572
$ws.Range("I86").Value = 1999.3846
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1999.3846
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -876.3846000000001
$ws.Range("H89").Value = 2070.8572
$ws.Range("I89").Value = 1999.3846
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 9996.923000000001
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -4380.923000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 3608.9092
$ws.Range("I15").Value = 600.75
$ws.Range("J15").Value = 5327.857
$ws.Range("K15").Value = 600.75
$ws.Range("L15").Value = 5327.857
$ws.Range("M15").Value = -430.75
$ws.Range("H31").Value = 21200.36
$ws.Range("I31").Value = 30272.914
$ws.Range("J31").Value = 3559.2778
$ws.Range("K31").Value = 30272.914
$ws.Range("L31").Value = 3559.2778
$ws.Range("M31").Value = -29977.914
$ws.Range("H34").Value = 21200.36
$ws.Range("I34").Value = 30272.914
$ws.Range("J34").Value = 3559.2778
$ws.Range("K34").Value = 30272.914
$ws.Range("L34").Value = 3559.2778
$ws.Range("M34").Value = -30070.914
$ws.Range("H45").Value = 21355.334
$ws.Range("I45").Value = 14067
$ws.Range("J45").Value = 24999.5
$ws.Range("K45").Value = 14067
$ws.Range("L45").Value = 24999.5
$ws.Range("M45").Value = -13474
$ws.Range("N45").Value = -26185.5
$ws.Range("H47").Value = 62500
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 62500
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 62500
$ws.Range("N47").Value = -63632
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H122").Value = 1952.9166
$ws.Range("I122").Value = 1867
$ws.Range("J122").Value = 2124.75
$ws.Range("K122").Value = 5601
$ws.Range("L122").Value = 6374.25
$ws.Range("M122").Value = -3151
$ws.Range("N122").Value = -11274.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 710.7646999999999
$ws.Range("I113").Value = 422
$ws.Range("J113").Value = 749.26666
$ws.Range("K113").Value = 1266
$ws.Range("L113").Value = 2247.79998
$ws.Range("M113").Value = 904
$ws.Range("N113").Value = -6587.79998
$ws.Range("H139").Value = 2562.75
$ws.Range("I139").Value = 2139.7778
$ws.Range("J139").Value = 3831.6667
$ws.Range("K139").Value = 6419.3334
$ws.Range("L139").Value = 11495.0001
$ws.Range("M139").Value = -1279.3334

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 982.5
$ws.Range("I31").Value = 982.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 982.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -690.5
$ws.Range("N31").ClearContents()
$ws.Range("H37").Value = 982.5
$ws.Range("I37").Value = 982.5
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 982.5
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -705.5
$ws.Range("N37").ClearContents()
$ws.Range("H52").Value = 15250
$ws.Range("I52").Value = 500
$ws.Range("J52").Value = 30000
$ws.Range("K52").Value = 500
$ws.Range("L52").Value = 30000
$ws.Range("M52").Value = -241
$ws.Range("N52").Value = -30518
$ws.Range("H113").Value = 1999.4445
$ws.Range("I113").Value = 1999.375
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1999.375
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 170.625
$ws.Range("H122").Value = 2288.5557
$ws.Range("I122").Value = 1571.8334
$ws.Range("J122").Value = 3722
$ws.Range("K122").Value = 4715.5002
$ws.Range("L122").Value = 11166
$ws.Range("M122").Value = -2265.5002
$ws.Range("H135").Value = 49635.816
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 49635.816
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 49635.816
$ws.Range("N135").Value = -59775.816

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4441.9
$ws.Range("I61").Value = 4386.2104
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 4386.2104
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -4184.2104
$ws.Range("H87").Value = 89811.42999999999
$ws.Range("I87").Value = 68500
$ws.Range("J87").Value = 93363.336
$ws.Range("K87").Value = 68500
$ws.Range("L87").Value = 93363.336
$ws.Range("M87").Value = -67377
$ws.Range("N87").Value = -95609.336
$ws.Range("H90").Value = 89811.42999999999
$ws.Range("I90").Value = 68500
$ws.Range("J90").Value = 93363.336
$ws.Range("K90").Value = 205500
$ws.Range("L90").Value = 280090.008
$ws.Range("M90").Value = -199884
$ws.Range("N90").Value = -291322.008
$ws.Range("H113").Value = 4441.9
$ws.Range("I113").Value = 4386.2104
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 4386.2104
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = -2216.2104
$ws.Range("H122").Value = 6476.8887
$ws.Range("I122").Value = 5448.8887
$ws.Range("J122").Value = 8532.888999999999
$ws.Range("K122").Value = 16346.6661
$ws.Range("L122").Value = 25598.667
$ws.Range("M122").Value = -13896.6661

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 31999
$ws.Range("I80").Value = 13000
$ws.Range("J80").Value = 35798.8
$ws.Range("K80").Value = 13000
$ws.Range("L80").Value = 35798.8
$ws.Range("M80").Value = -12002
$ws.Range("N80").Value = -37794.8
$ws.Range("H81").Value = 8627.817999999999
$ws.Range("I81").Value = 21588.4
$ws.Range("J81").Value = 4815.8823
$ws.Range("K81").Value = 43176.8
$ws.Range("L81").Value = 9631.7646
$ws.Range("M81").Value = -42115.8
$ws.Range("N81").Value = -11753.7646
$ws.Range("H83").Value = 31999
$ws.Range("I83").Value = 13000
$ws.Range("J83").Value = 35798.8
$ws.Range("K83").Value = 39000
$ws.Range("L83").Value = 107396.4
$ws.Range("M83").Value = -34008
$ws.Range("N83").Value = -117380.4
$ws.Range("H84").Value = 8627.817999999999
$ws.Range("I84").Value = 21588.4
$ws.Range("J84").Value = 4815.8823
$ws.Range("K84").Value = 215884
$ws.Range("L84").Value = 48158.823
$ws.Range("M84").Value = -210580
$ws.Range("N84").Value = -58766.823
$ws.Range("H122").Value = 2731.3333
$ws.Range("I122").Value = 2543
$ws.Range("J122").Value = 3899
$ws.Range("K122").Value = 7629
$ws.Range("L122").Value = 11697
$ws.Range("M122").Value = -5179
$ws.Range("H136").Value = 1813.5834
$ws.Range("I136").Value = 1458.2258
$ws.Range("J136").Value = 2461.5881
$ws.Range("K136").Value = 4374.6774
$ws.Range("L136").Value = 7384.7643
$ws.Range("M136").Value = -1824.6774
$ws.Range("H137").Value = 100683
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 100683
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 100683
$ws.Range("N137").Value = -110883
